$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 2
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = ";0;0"
$ws.Range("G9").Value = ";13;23"
$ws.Range("H9").Value = ";-1000;-98000.0"
